# Update the "Fitness" values (column C) for rows 2-110 with the new
# best-fitness-so-far figures recorded for this run.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newFitness = @(
    12201,12090,9853,9561,9561,9561,9561,9561,9352,9352,
    9352,9162,9162,9162,8779,8779,8779,8779,8779,8779,
    8581,8581,8581,8581,8581,8548,8318,7639,7639,7639,
    7639,7639,7639,7639,7639,7639,7639,7639,7639,7639,
    7639,7639,7639,7639,7639,7639,7639,7639,7639,7639,
    7639,7639,7639,7639,7639,7639,7639,7639,7639,7639,
    7639,7639,7639,7639,7639,7581,7581,7581,7581,7581,
    7581,7534,7534,7534,7534,7534,7534,7534,7534,7534,
    7534,7534,7534,7534,7534,7534,7295,7295,7295,7295,
    7295,7295,7295,7295,7295,7295,7295,7295,7295,7295,
    7295,7295,7295,7295,7295,7295,7295,7295,7295
)

$startRow = 2
for ($i = 0; $i -lt $newFitness.Length; $i++) {
    $ws.Cells.Item($startRow + $i, 3).Value = $newFitness[$i]
}
